# Refresh the cryptocurrency price / 1h-volume table with the latest scrape
# (mirrors the GitHub Actions data-refresh job run on 2023-04-30).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "30.108.55"
$ws.Cells.Item(2, 5).Value = "  +2.01%  "

$ws.Cells.Item(3, 4).Value = "1.949.04"
$ws.Cells.Item(3, 5).Value = "  +1.31%  "

$ws.Cells.Item(4, 5).Value = "  +0.14%  "

$ws.Cells.Item(5, 4).Value = "'327.81"
$ws.Cells.Item(5, 5).Value = "  +0.66%  "

$ws.Cells.Item(6, 5).Value = "  +0.17%  "

$ws.Cells.Item(7, 5).Value = "  +0.30%  "

$ws.Cells.Item(8, 4).Value = "'0.4120"
$ws.Cells.Item(8, 5).Value = "  +0.48%  "

$ws.Cells.Item(9, 4).Value = "'0.08229"
$ws.Cells.Item(9, 5).Value = "  +0.48%  "

$ws.Cells.Item(10, 4).Value = "'1.022"

$ws.Cells.Item(11, 4).Value = "'24.08"
$ws.Cells.Item(11, 5).Value = "  +2.09%  "

$ws.Cells.Item(12, 4).Value = "1.956.74"
$ws.Cells.Item(12, 5).Value = "  +2.90%  "

$ws.Cells.Item(13, 4).Value = "'6.120"
$ws.Cells.Item(13, 5).Value = "  +1.04%  "

$ws.Cells.Item(14, 4).Value = "'7.348"
$ws.Cells.Item(14, 5).Value = "  +1.23%  "

$ws.Cells.Item(15, 4).Value = "'92.03"
$ws.Cells.Item(15, 5).Value = "  +0.58%  "

$ws.Cells.Item(16, 4).Value = "'0.06858"
$ws.Cells.Item(16, 5).Value = "  +1.23%  "

$ws.Cells.Item(17, 4).Value = "'1.010"
$ws.Cells.Item(17, 5).Value = "  +0.18%  "

$ws.Cells.Item(18, 5).Value = "  -0.04%  "

$ws.Cells.Item(19, 4).Value = "'18.00"
$ws.Cells.Item(19, 5).Value = "  +0.96%  "

$ws.Cells.Item(20, 4).Value = "'1.008"
$ws.Cells.Item(20, 5).Value = "  +0.20%  "

$ws.Cells.Item(21, 4).Value = "30.115.34"
$ws.Cells.Item(21, 5).Value = "  +1.97%  "

$ws.Cells.Item(22, 4).Value = "'5.689"
$ws.Cells.Item(22, 5).Value = "  +0.88%  "

$ws.Cells.Item(23, 4).Value = "'11.99"
$ws.Cells.Item(23, 5).Value = "  +1.59%  "

$ws.Cells.Item(24, 5).Value = "  +0.45%  "

$ws.Cells.Item(25, 4).Value = "2.173.32"
$ws.Cells.Item(25, 5).Value = "  +1.86%  "

$ws.Cells.Item(26, 2).Value = "Monero"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(26, 4).Value = "'156.94"
$ws.Cells.Item(26, 5).Value = "  +0.16%  "

$ws.Cells.Item(27, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(27, 4).Value = "'6.534"
$ws.Cells.Item(27, 5).Value = "  -3.54%  "

$ws.Cells.Item(28, 4).Value = "'20.17"
$ws.Cells.Item(28, 5).Value = "  +0.29%  "

$ws.Cells.Item(29, 4).Value = "'2.121"
$ws.Cells.Item(29, 5).Value = "  -0.22%  "

$ws.Cells.Item(30, 4).Value = "'121.42"
$ws.Cells.Item(30, 5).Value = "  +0.64%  "

$ws.Cells.Item(31, 4).Value = "'1.028"
$ws.Cells.Item(31, 5).Value = "  -0.43%  "

$ws.Cells.Item(32, 4).Value = "'0.09653"
$ws.Cells.Item(32, 5).Value = "  +0.64%  "

$ws.Cells.Item(33, 4).Value = "'5.657"
$ws.Cells.Item(33, 5).Value = "  +2.14%  "

$ws.Cells.Item(34, 4).Value = "'1.433"
$ws.Cells.Item(34, 5).Value = "  +2.85%  "

$ws.Cells.Item(35, 4).Value = "'3.563"
$ws.Cells.Item(35, 5).Value = "  -0.19%  "

$ws.Cells.Item(36, 4).Value = "'0.06529"
$ws.Cells.Item(36, 5).Value = "  +6.15%  "

$ws.Cells.Item(37, 4).Value = "'0.02308"
$ws.Cells.Item(37, 5).Value = "  +0.90%  "

$ws.Cells.Item(38, 4).Value = "'1.230"
$ws.Cells.Item(38, 5).Value = "  +3.89%  "

$ws.Cells.Item(39, 4).Value = "'0.5986"
$ws.Cells.Item(39, 5).Value = "  -0.15%  "

$ws.Cells.Item(40, 4).Value = "'10.81"
$ws.Cells.Item(40, 5).Value = "  -0.02%  "

$ws.Cells.Item(41, 4).Value = "'8.000"
$ws.Cells.Item(41, 5).Value = "  -0.71%  "

$ws.Cells.Item(42, 4).Value = "'2.538"
$ws.Cells.Item(42, 5).Value = "  +5.49%  "

$ws.Cells.Item(43, 4).Value = "'0.1862"
$ws.Cells.Item(43, 5).Value = "  -0.21%  "

$ws.Cells.Item(44, 4).Value = "'1.283"
$ws.Cells.Item(44, 5).Value = "  -0.07%  "

$ws.Cells.Item(45, 4).Value = "'12.47"
$ws.Cells.Item(45, 5).Value = "  +0.30%  "

$ws.Cells.Item(46, 4).Value = "'0.07567"
$ws.Cells.Item(46, 5).Value = "  -0.61%  "

$ws.Cells.Item(47, 4).Value = "'0.5599"
$ws.Cells.Item(47, 5).Value = "  -0.03%  "

$ws.Cells.Item(48, 4).Value = "'2.000"
$ws.Cells.Item(48, 5).Value = "  +1.87%  "

$ws.Cells.Item(49, 4).Value = "'118.29"
$ws.Cells.Item(49, 5).Value = "  +0.61%  "

$ws.Cells.Item(50, 4).Value = "'2.441"
$ws.Cells.Item(50, 5).Value = "  +0.00%  "

$ws.Cells.Item(51, 4).Value = "'72.84"
$ws.Cells.Item(51, 5).Value = "  +0.00%  "

Write-Host "Applied cryptos update: 50 rows refreshed"